$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $ok = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    Write-Output "Replace [$old] -> [$new] : $ok"
}

# Title heading (also appears again later as a bold run near the end of the doc)
Replace-Text "Play Dance Party Free: Review and Gameplay" "Play Dance Party Free - Review"

# "What we like" bullet list
Replace-Text "Free Spins with Progressive Multiplier up to 30x" "Scatter symbol and free spins bonus"
Replace-Text "243 ways to win in 3x5 format" "243 ways to win"
Replace-Text "Bet range from 0.01-0.50 cents up to €100" "Betting options"
Replace-Text "Smooth and energetic gaming experience" "RTP and gaming experience"

# "What we don't like" bullet list
Replace-Text "Static animations" "Limited symbols and animations"
Replace-Text "Lack of bonus features aside from free spins" "Structure remains old school"

# Closing italic summary paragraph
Replace-Text "Read our review of Dance Party, an online slot game with free spins and 243 ways to win. Play for free and learn about this energetic, smooth game." "Read our review of Dance Party and play for free. Discover the scatter symbol, free spins, and more."
